$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from the existing header cell (H1) so the new headers match formatting
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for rows 2-16 in columns I (I0) and J (IF)
$data = @(
    @(10, 10),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(5, 6),
    @(6, 6),
    @(5, 5),
    @(8, 8),
    @(5, 5),
    @(9, 9),
    @(9, 9),
    @(5, 5),
    @(5, 6),
    @(4, 4),
    @(3, 3)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
